# Apply updated weekly Fruta/Hortaliza data (Femacal de La Calera - Perejil)
# Columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) are reshuffled across
# rows 2-24 (the rest of the columns stay identical between rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @(D, J, K, L, M, P)
$data = @{
    2  = @(44340, 54,  3000, 3000, 3000, 1000)
    3  = @(44389, 81,  2800, 3000, 2889, 963)
    4  = @(44243, 45,  3000, 3000, 3000, 1000)
    5  = @(44223, 80,  2500, 3000, 2781, 927)
    6  = @(44291, 45,  3000, 3000, 3000, 1000)
    7  = @(44292, 40,  3000, 3000, 3000, 1000)
    8  = @(44537, 88,  2000, 2200, 2091, 697)
    9  = @(44179, 78,  3000, 3000, 3000, 1000)
    10 = @(44222, 45,  3000, 3000, 3000, 1000)
    11 = @(44559, 68,  2000, 2000, 2000, 667)
    12 = @(44224, 67,  3000, 3000, 3000, 1000)
    13 = @(44390, 50,  3000, 3000, 3000, 1000)
    14 = @(44574, 50,  3000, 3000, 3000, 1000)
    15 = @(44165, 68,  3000, 3000, 3000, 1000)
    16 = @(44557, 104, 2000, 2500, 2260, 753)
    17 = @(44242, 95,  2500, 3000, 2737, 912)
    18 = @(44193, 70,  3000, 3000, 3000, 1000)
    19 = @(44166, 45,  2500, 2500, 2500, 833)
    20 = @(44536, 125, 2200, 2200, 2200, 733)
    21 = @(44187, 65,  3000, 3000, 3000, 1000)
    22 = @(44225, 56,  3000, 3000, 3000, 1000)
    23 = @(44221, 50,  2500, 2500, 2500, 833)
    24 = @(44260, 60,  3500, 3500, 3500, 1167)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $ws.Cells.Item($row, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]   # P - Precio $/Kg
}
